# Updates crypto price/volume figures to match the latest GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.253.59"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.12"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.57%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.67%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2865"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06557"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.17"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +10.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07904"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "98.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.868.56"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.183"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6814"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "277.95"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.251.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.62"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +8.11%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.116.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.365"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.197"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.234"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.07"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.951"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09847"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.384"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.482"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.071"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04749"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.138"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7045"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.708"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01880"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.630"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "75.77"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.288"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.956"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8555"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4176"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.0000"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.41"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.228"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "947.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.234"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.25"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.12%  "
